$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the series: insert a fresh
# row at position 184 (pushing the existing rows 184-220 down to 185-221)
# and populate it with the new record.
$ws.Rows.Item(184).Insert()

$ws.Cells.Item(184, 1).Value = 3
$ws.Cells.Item(184, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44476
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 6).Value = 100112031
$ws.Cells.Item(184, 7).Value = "Poroto verde"
$ws.Cells.Item(184, 8).Value = "Magnum"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 35
$ws.Cells.Item(184, 11).Value = 38000
$ws.Cells.Item(184, 12).Value = 38000
$ws.Cells.Item(184, 13).Value = 38000
$ws.Cells.Item(184, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(184, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(184, 16).Value = 1520
$ws.Cells.Item(184, 17).Value = 25
$ws.Cells.Item(184, 18).Value = "Hortaliza"
